# Insert a new row for "Port-au-Prince, Haiti" (colo PAP) above the current
# row 223 (Amman, Jordan), shifting Amman and everything below it down by one
# row. This mirrors the source data gaining a new data-center entry that
# sorts alphabetically before Amman.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 223; existing row 223 (and all following
# rows) shift down to row 224, etc.
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with the Port-au-Prince, Haiti data.
$ws.Cells.Item(223, 1).Value = "PAP"
$ws.Cells.Item(223, 2).Value = "Port-au-Prince, Haiti"
$ws.Cells.Item(223, 3).Value = 18.5799999237
$ws.Cells.Item(223, 4).Value = -72.2925033569
$ws.Cells.Item(223, 5).Value = "HT"
$ws.Cells.Item(223, 6).Value = "North America"
$ws.Cells.Item(223, 7).Value = "Port-au-Prince"

# Match the bold/bordered/centered style used by the rest of column A
# (font bold, thin box border, centered horizontally, top-aligned vertically).
$colAHeaderCell = $ws.Cells.Item(223, 1)
$colAHeaderCell.Font.Bold = $true
$colAHeaderCell.HorizontalAlignment = -4108  # xlCenter
$colAHeaderCell.VerticalAlignment = -4160    # xlTop
$colAHeaderCell.Borders.LineStyle = 1        # xlContinuous
$colAHeaderCell.Borders.Weight = 2           # xlThin
